# Add the new "2022-Q3" worksheet and data (as described by the commit
# "feat: add 2022-Q3 data"), shifting the existing quarterly sheets along.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new sheet right after "总计", before "2022-Q1".
#    (Worksheet references are (re-)fetched *after* Add(), since adding a
#    sheet invalidates handles obtained beforehand.)
# ---------------------------------------------------------------------
$totalTmp = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($null, $totalTmp)
$newSheet.Name = "2022-Q3"

$total = $wb.Worksheets.Item("总计")
$q1_2022 = $wb.Worksheets.Item("2022-Q1")
$newSheet = $wb.Worksheets.Item("2022-Q3")

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: existing rows 2-7 shift down to
#    3-8, and a brand new row 2 is written for 2022-Q3.
# ---------------------------------------------------------------------
$summaryDates = @("2022-Q3", "2022-Q1", "2021-Q4", "2021-Q3", "2021-Q2", "2021-Q1", "2020-Q4")
$summaryCounts = @(7, 4, 6, 5, 9, 11, 6)
$summaryValues = @(0.09, 0.07, 0.05, 0.07, 0.32, 0.51, 0.05)

for ($i = 0; $i -lt 7; $i++) {
    $r = $i + 2
    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 2).Value = $summaryDates[$i]
    $total.Cells.Item($r, 3).Value = $summaryCounts[$i]
    $total.Cells.Item($r, 4).Value = $summaryValues[$i]
}
# Row 8 (2020-Q4) is a brand-new row - give its index cell (A8) the same
# style as the other index cells in column A.
$total.Range("A7").Copy($total.Range("A8"))
$total.Cells.Item(8, 1).Value = 6

# ---------------------------------------------------------------------
# 3. Populate the new "2022-Q3" sheet with the fund-holding table, using
#    the "2022-Q1" sheet's header/index-column formatting as a template.
# ---------------------------------------------------------------------
$q1_2022.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$q1_2022.Range("A2").Copy($newSheet.Range("A2:A8"))

$codes = @("501305", "513530", "501306", "006658", "501307", "006659", "501308")
$names = @(
    "汇添富中证港股通高股息投资指数（LOF）A",
    "华泰柏瑞中证港股通高股息投资ETF（QDII）",
    "汇添富中证港股通高股息投资指数（LOF）C",
    "财通中证香港红利等权投资指数A",
    "银河中证沪港深高股息指数（LOF）A",
    "财通中证香港红利等权投资指数C",
    "银河中证沪港深高股息指数（LOF）C"
)
$scales = @("0.87", "0.78", "0.17", "0.13", "0.15", "0.04", "0.01")
$stockPos = @("92.21", "95.80", "92.21", "88.79", "90.33", "88.79", "90.33")
$posRatio = @("4.32", "4.62", "4.32", "3.68", "1.41", "3.68", "1.41")
$heldValue = @("0.0376", "0.0360", "0.0073", "0.0048", "0.0021", "0.0015", "0.0001")
$rank = @(3, 3, 3, 7, 6, 7, 6)

# Force columns B:G to be stored as text (matches the source file, where
# fund codes / percentages are inline strings, not numbers) without
# leaving a lingering explicit number format on the cells.
$textRange = $newSheet.Range("B2:G8")
$textRange.NumberFormat = "@"

for ($i = 0; $i -lt 7; $i++) {
    $r = $i + 2
    $newSheet.Cells.Item($r, 1).Value = $i
    $newSheet.Cells.Item($r, 2).Value = $codes[$i]
    $newSheet.Cells.Item($r, 3).Value = $names[$i]
    $newSheet.Cells.Item($r, 4).Value = $scales[$i]
    $newSheet.Cells.Item($r, 5).Value = $stockPos[$i]
    $newSheet.Cells.Item($r, 6).Value = $posRatio[$i]
    $newSheet.Cells.Item($r, 7).Value = $heldValue[$i]
    $newSheet.Cells.Item($r, 8).Value = $rank[$i]
}

# Drop the temporary "@" number format back to the default style now that
# the text has been committed, so the cells stay styleless like the rest
# of the workbook.
$textRange.Style = "Normal"

# ---------------------------------------------------------------------
# 4. Keep "2020-Q4" as the active/selected tab, matching the source file.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
